$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19; this shifts existing rows 19:116 down to 20:117
# and copies formatting (incl. date style on column D) from the row above.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new data record.
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value = 44670
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 100112032
$ws.Range("G19").Value = "Zapallo italiano"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 11000
$ws.Range("M19").Value = 10444
$ws.Range("N19").Value = "$/caja 50 unidades"
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 209
$ws.Range("Q19").Value = 50
$ws.Range("R19").Value = "Hortaliza"
